# Append 45 new rows (102-146) to the master-reg_center_machine_device sheet,
# following the same pattern as the existing rows, then update the sheet's
# active view (scroll/selection) to match where the user left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(102, 10002, 10021, 3000121),
    @(103, 10003, 10022, 3000122),
    @(104, 10004, 10023, 3000123),
    @(105, 10005, 10024, 3000124),
    @(106, 10006, 10025, 3000125),
    @(107, 10007, 10026, 3000126),
    @(108, 10008, 10027, 3000127),
    @(109, 10009, 10028, 3000128),
    @(110, 10010, 10029, 3000129),
    @(111, 10002, 10021, 3000130),
    @(112, 10003, 10022, 3000131),
    @(113, 10004, 10023, 3000132),
    @(114, 10005, 10024, 3000133),
    @(115, 10006, 10025, 3000134),
    @(116, 10007, 10026, 3000135),
    @(117, 10008, 10027, 3000136),
    @(118, 10009, 10028, 3000137),
    @(119, 10010, 10029, 3000138),
    @(120, 10002, 10021, 3000139),
    @(121, 10003, 10022, 3000140),
    @(122, 10004, 10023, 3000141),
    @(123, 10005, 10024, 3000142),
    @(124, 10006, 10025, 3000143),
    @(125, 10007, 10026, 3000144),
    @(126, 10008, 10027, 3000145),
    @(127, 10009, 10028, 3000146),
    @(128, 10010, 10029, 3000147),
    @(129, 10002, 10021, 3000148),
    @(130, 10003, 10022, 3000149),
    @(131, 10004, 10023, 3000150),
    @(132, 10005, 10024, 3000151),
    @(133, 10006, 10025, 3000152),
    @(134, 10007, 10026, 3000153),
    @(135, 10008, 10027, 3000154),
    @(136, 10009, 10028, 3000155),
    @(137, 10010, 10029, 3000156),
    @(138, 10002, 10021, 3000157),
    @(139, 10003, 10022, 3000158),
    @(140, 10004, 10023, 3000159),
    @(141, 10005, 10024, 3000160),
    @(142, 10006, 10025, 3000161),
    @(143, 10007, 10026, 3000162),
    @(144, 10008, 10027, 3000163),
    @(145, 10009, 10028, 3000164),
    @(146, 10010, 10029, 3000165),
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = "eng"
    $ws.Cells.Item($rowNum, 5).Value = $true
    $ws.Cells.Item($rowNum, 6).Value = "superadmin()"
    $ws.Cells.Item($rowNum, 7).Value = "now()"
}

# Match the saved view state: selection spanning the newly added rows,
# scrolled so row 129 is at the top of the window.
$ws.Range("A102:G146").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 129

# Printer/page setup metadata recorded on save.
$ws.PageSetup.Orientation = 1
